$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 - new entry for 6-Mar-2021 (serial 44261)
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 1).Interior.Pattern = -4142

$ws.Cells.Item(41, 2).Value = 26426379

$ws.Cells.Item(41, 3).Value = 36956
$ws.Cells.Item(41, 3).NumberFormat = "#,##0.00"

$ws.Cells.Item(41, 4).Value = 34999.22
$ws.Cells.Item(41, 4).NumberFormat = "#,##0.00"

$ws.Cells.Item(41, 5).Value = 44261
$ws.Cells.Item(41, 5).NumberFormat = "[`$-409]dd\-mmm\-yy;@"

# Row 42 - new entry for 7-Mar-2021 (serial 44262)
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 1).Interior.Pattern = -4142

$ws.Cells.Item(42, 2).Value = 26446577

$ws.Cells.Item(42, 3).Value = 36956
$ws.Cells.Item(42, 3).NumberFormat = "#,##0.00"

$ws.Cells.Item(42, 4).Value = 34999.22
$ws.Cells.Item(42, 4).NumberFormat = "#,##0.00"

$ws.Cells.Item(42, 5).Value = 44262
$ws.Cells.Item(42, 5).NumberFormat = "[`$-409]dd\-mmm\-yy;@"

# Re-enter the carried-down formulas so the previously-blank shared-formula
# cells (which cached an empty string result) recompute against the new data
$ws.Cells.Item(41, 6).Formula = '=IF(B41="","",C41-D41)'
$ws.Cells.Item(41, 7).Formula = '=IF(B41="","",F41/D41*100)'
$ws.Cells.Item(41, 8).Formula = '=IF(B41="","",D41*1.04)'
$ws.Cells.Item(41, 9).Formula = '=IF(B41="","",C41-H41)'

$ws.Cells.Item(42, 6).Formula = '=IF(B42="","",C42-D42)'
$ws.Cells.Item(42, 7).Formula = '=IF(B42="","",F42/D42*100)'
$ws.Cells.Item(42, 8).Formula = '=IF(B42="","",D42*1.04)'
$ws.Cells.Item(42, 9).Formula = '=IF(B42="","",C42-H42)'

# Move the active selection to where the author left off
$ws.Range("J42").Select()
